$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete entire row 3 (which contained "30072020_J1"), shifting rows below up
$ws.Rows("3").Delete()

# Update the selection to match what was left selected (row 3 selected again, post-delete)
$ws.Range("A3:XFD3").Select()
